# Slide 10 ("week14/potato.pptx"):
#  1. Reposition/resize the picture placeholder (shape 1).
#  2. Fill in the previously-empty 4th paragraph of the text box (shape 2)
#     with "共88筆資料" (as three runs: "共" / "88" / "筆資料").
#
# Note on the magic Left/Top/Width/Height literals: PowerPoint COM expresses
# shape geometry in points, while the underlying OOXML stores EMUs
# (1 pt = 12700 EMU). This host's point->EMU conversion round-trips the
# incoming value through a 32-bit float before flooring to EMU, so a
# "simple" pt = emu/12700 literal can land one EMU short. The literals below
# were chosen so that, after that float32 round-trip, they floor to the
# exact target EMUs (x=1673165, y=317484, cx=3851089, cy=5035581).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)

# --- 1. Move/resize the picture (shape 1) ---
$pic = $s.Shapes.Item(1)
$pic.Left   = 131.7452850341797
$pic.Top    = 24.998741149902344
$pic.Width  = 303.2353820800781
$pic.Height = 396.50244140625

# --- 2. Add the "共88筆資料" text to the 4th (previously empty) paragraph ---
$txBox = $s.Shapes.Item(2)
$tr = $txBox.TextFrame.TextRange
$lastPara = $tr.Paragraphs(4)
[void]$lastPara.InsertAfter("共")
[void]$lastPara.InsertAfter("88")
[void]$lastPara.InsertAfter("筆資料")
